$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 18
$ws.Range("I18").Value = "sv"
$ws.Range("J18").Value = "Statement-opinion"

# Row 21
$ws.Range("I21").Value = "sv"
$ws.Range("J21").Value = "Statement-opinion"

# Row 41
$ws.Range("I41").Value = "sd"
$ws.Range("J41").Value = "Statement-non-opinion"

# Row 45
$ws.Range("I45").Value = "sd"
$ws.Range("J45").Value = "Statement-non-opinion"

# Row 52
$ws.Range("I52").Value = "sv"
$ws.Range("J52").Value = "Statement-opinion"

# Row 62
$ws.Range("I62").Value = "sd"
$ws.Range("J62").Value = "Statement-non-opinion"
